$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the definition in C4
$ws.Range("C4").Value = "To spend time with one's loved ones"

# Remove row 5 (fug / fugg / long definition) entirely
$ws.Rows("5:5").Delete()

# Match the final selection state left by the editing user
$ws.Range("C4").Select()
